$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 320, shifting the existing rows 320:343
# down to 321:344 (preserves their data/styles, extends used range to R344).
$ws.Rows.Item(320).Insert()

# Populate the newly inserted row 320 with the new weekly record.
$ws.Cells.Item(320,1).Value  = 6
$ws.Cells.Item(320,2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(320,3).Value  = "Metropolitana"
$ws.Cells.Item(320,4).Value  = 45166
$ws.Cells.Item(320,5).Value  = 13
$ws.Cells.Item(320,6).Value  = 100112022
$ws.Cells.Item(320,7).Value  = "Arveja Verde"
$ws.Cells.Item(320,8).Value  = "Perfection"
$ws.Cells.Item(320,9).Value  = "Primera"
$ws.Cells.Item(320,10).Value = 400
$ws.Cells.Item(320,11).Value = 25000
$ws.Cells.Item(320,12).Value = 27000
$ws.Cells.Item(320,13).Value = 25850
$ws.Cells.Item(320,14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(320,15).Value = "Región Metropolitana"
$ws.Cells.Item(320,16).Value = 1034
$ws.Cells.Item(320,17).Value = 25
$ws.Cells.Item(320,18).Value = "Hortaliza"
